$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 337; this shifts existing rows 337-355 down to 338-356.
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new data record.
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = "Femacal de La Calera"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value2 = 44706
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100114013
$ws.Range("G337").Value = "Zanahoria"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 340
$ws.Range("K337").Value = 7000
$ws.Range("L337").Value = 7500
$ws.Range("M337").Value = 7279
$ws.Range("N337").Value = "$/saco 20 kilos"
$ws.Range("O337").Value = "Provincia de Quillota"
$ws.Range("P337").Value = 364
$ws.Range("Q337").Value = 20
$ws.Range("R337").Value = "Hortaliza"
